# guardduty-org.pptx edit script
# Commit: "Add fill to allow image to be visible in dark mode"
#  - Re-cache the "12/5/21" datetimeFigureOut date placeholders (slide master +
#    every layout + notes master) that PowerPoint stamps on every save.
#  - Give the big full-bleed frame rectangle behind the architecture diagram a
#    solid bg1 fill (was noFill) so the picture is visible in dark mode.
#  - Rename the "OU: Core" label to "OU: Security" on the two account boxes.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the auto date field text ("11/16/21" -> "12/5/21") everywhere a
#    "Date Placeholder *" shape appears: slide master, each slide layout, and
#    the notes master.
# ---------------------------------------------------------------------------
function Set-DatePlaceholderText($shapes, $newText) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

$newDate = "12/5/21"

$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes $newDate

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Set-DatePlaceholderText $layouts.Item($i).Shapes $newDate
}

Set-DatePlaceholderText $p.NotesMaster.Shapes $newDate

# ---------------------------------------------------------------------------
# 2) Give "Rectangle 124" (the big frame behind the architecture diagram on
#    slide 1) a solid bg1 fill instead of no fill.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$frame = $slide.Shapes.Item("Rectangle 124")
$frame.Fill.ForeColor.SchemeColor = "bg1"

# ---------------------------------------------------------------------------
# 3) Rename "OU: Core" -> "OU: Security" on the two account rectangles.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -like "*OU: Core*") {
            $paragraphs = $tr.Paragraphs()
            for ($k = 1; $k -le $paragraphs.Count; $k++) {
                $para = $tr.Paragraphs($k)
                if ($para.Text -eq "OU: Core") {
                    $run = $para.Runs(1)
                    $run.Text = "OU: Security"
                }
            }
        }
    }
}
